# Automatic update of files.
# The edit re-shuffles the per-observation data among rows 19-25 and 27-30
# of the "Artfynd" sheet (row 26 is untouched). Only the columns that
# actually carry observation-specific data are affected:
#   A (Id), B (Taxonsorteringsordning), D (Rödlistade), E (TaxonId),
#   F (Artnamn), G (Vetenskapligt namn), H (Auktor), Q (Ost), R (Nord)
# Every other column is identical across these rows, so it is left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by 1-based index) that carry the row-specific content.
$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18)

# Snapshot the current ("before") values for every affected row so that the
# re-shuffle below can be performed without clobbering data that is still
# needed for a later assignment.
$rowsToRead = @(19, 20, 21, 22, 23, 24, 25, 27, 28, 29, 30)
$snapshot = @{}
foreach ($r in $rowsToRead) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Mapping of destination row -> source row (i.e. the content that used to
# live in the source row now belongs to the destination row).
$rowMap = @{
    19 = 22
    20 = 21
    21 = 28
    22 = 19
    23 = 20
    24 = 25
    25 = 24
    27 = 29
    28 = 23
    29 = 30
    30 = 27
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}

"Row data re-shuffled successfully."
